# Apply edits from commit "Upload excel files with prices":
#  - Update the scrape timestamp (column O) for every data row (2..410)
#    from "2022-12-23 12:55:51" to "2022-12-23 20:49:23"
#  - Rename two Roland cracker products (title/href/aria-label) in rows 352 and 397

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2022-12-23 12:55:51"
$newTimestamp = "2022-12-23 20:49:23"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 15)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}

# Row 352: "Roland Cracker Sport High Protein Choco" -> "Roland Sport High Protein Choco"
$ws.Range("B352").Value = "Roland Sport High Protein Choco"
$ws.Range("C352").Value = "/de/lebensmittel/brot-backwaren/cracker-knaeckebrot/cracker/suess/roland-sport-high-protein-choco/p/6729458"
$ws.Range("M352").Value = "Roland Sport High Protein Choco 4.60 Schweizer Franken"

# Row 397: "Roland Cracker Protein Nature" -> "Roland Sport Plus Nature"
$ws.Range("B397").Value = "Roland Sport Plus Nature"
$ws.Range("C397").Value = "/de/lebensmittel/brot-backwaren/cracker-knaeckebrot/cracker/nature/roland-sport-plus-nature/p/6729401"
$ws.Range("M397").Value = "Roland Sport Plus Nature 4.60 Schweizer Franken"
